# Simulate three more form submissions landing in the "Captured_Values"
# sheet: each one is a new row with the captured number (A) and the
# captured text (B), appended right after the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

$capturedNumber = 123456789
$capturedText   = "Real Programmers Count 0123456789 From Zero"

# Existing data runs from row 2 through the last populated row; find it
# the same way a real macro would (Ctrl+Up from the bottom of column A).
$lastRow  = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt 3; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $capturedNumber
    $ws.Cells.Item($row, 2).Value = $capturedText
}
